# Fixed LHS sampling to only sample across uncertainties (X) that vary (Ls still vary for all Ls)
# and rebuilt templates with PFLO:ALL_NO_STOPPING_DEFORESTATION_PLUR
#
# This script:
#  1. Renames the existing "strategy_id-5008" worksheet to "strategy_id-5007"
#  2. Adds a new worksheet named "strategy_id-5009" right after it, containing
#     the same data as "strategy_id-5007" (a duplicate of the strategy template).

$wb = $excel.ActiveWorkbook

# 1) Rename strategy_id-5008 -> strategy_id-5007
$srcSheet = $wb.Worksheets.Item("strategy_id-5008")
$srcSheet.Name = "strategy_id-5007"

# 2) Duplicate the sheet, placing the copy immediately after it, then rename
#    the new copy to strategy_id-5009
$srcSheet.Copy($null, $srcSheet)
$newSheet = $wb.Worksheets.Item($srcSheet.Index + 1)
$newSheet.Name = "strategy_id-5009"
